$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*System*") {
        $parts = $val -split ",\s*"
        $sysIndex = -1
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($parts[$i].Equals("System")) {
                $sysIndex = $i
                break
            }
        }
        if ($sysIndex -gt 0) {
            $newVal = $parts[$sysIndex]
            for ($i = 0; $i -lt $parts.Length; $i++) {
                if ($i -ne $sysIndex) {
                    $newVal = $newVal + ", " + $parts[$i]
                }
            }
            $cell.Value2 = $newVal
        }
    }
}
